$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L2").Value = "*maa://24633 (56.05), *maa://30515 (69.31), *maa://34787 (72.86), ***maa://20792 (11.93), maa://39402 (88.64), ***maa://29083 (27.78)"
$ws.Range("AF2").Value = "maa://25251 (91.84), ***maa://21730 (22.86), ***maa://39501 (18.18), *maa://36675 (60.0)"
$ws.Range("L3").Value = "*maa://22880 (65.95), maa://20276 (85.09), *maa://22749 (72.73)"
$ws.Range("D7").Value = "maa://21955 (94.29)"
$ws.Range("H7").Value = "*maa://22763 (66.67)"
$ws.Range("A8").Value = "更新日期：2025.01.01 13:19:57"
$ws.Range("D8").Value = "*maa://21476 (72.92), **maa://39431 (45.45), *maa://37551 (57.14)"
$ws.Range("AB9").Value = "maa://28711 (87.38), ***maa://22740 (5.77), **maa://39938 (48.0), **maa://27377 (42.86), ***maa://25174 (19.05), maa://40166 (93.33)"
$ws.Range("T10").Value = "maa://27395 (96.09), maa://22755 (87.61), **maa://22756 (40.91), ***maa://21737 (10.61)"
$ws.Range("X10").Value = "maa://22301 (97.64), maa://22726 (100.0)"
$ws.Range("L14").Value = "maa://26245 (96.53), maa://21288 (96.3), maa://39841 (95.18), maa://36682 (97.37)"
$ws.Range("D15").Value = "*maa://22743 (77.44), maa://22734 (84.03), *maa://30808 (65.08), **maa://36048 (32.61)"
$ws.Range("AF15").Value = "maa://21364 (80.91), *maa://22766 (70.91), *maa://36666 (78.82)"
$ws.Range("T16").Value = "maa://22729 (94.77), *maa://28648 (68.85), maa://36674 (82.93)"
$ws.Range("H18").Value = "maa://24421 (90.12)"
$ws.Range("AB19").Value = "*maa://30709 (63.59), *maa://36668 (55.84)"
$ws.Range("L23").Value = "maa://39756 (94.32), maa://39875 (93.75)"
$ws.Range("D24").Value = "*maa://24368 (78.89)"
$ws.Range("X24").Value = "maa://29988 (86.84), maa://23504 (93.15), **maa://22892 (39.31), *maa://25141 (76.38), *maa://36663 (78.26), ***maa://22815 (23.08)"
$ws.Range("X25").Value = "*maa://29890 (77.27)"
$ws.Range("AB25").Value = "maa://31215 (85.86), *maa://24516 (79.78), maa://26001 (87.5)"
$ws.Range("X26").Value = "maa://24389 (96.43)"
$ws.Range("H27").Value = "**maa://21283 (48.0), maa://34494 (96.55), *maa://39601 (76.47), **maa://36665 (44.44)"
$ws.Range("X28").Value = "maa://39929 (89.72), ***maa://39723 (14.29), maa://41749 (91.53)"
$ws.Range("AF28").Value = "maa://36660 (93.02), *maa://36701 (64.29)"
$ws.Range("L29").Value = "maa://28432 (92.93), *maa://28440 (76.84), maa://31400 (100.0), *maa://28650 (71.43)"
$ws.Range("P29").Value = "*maa://23168 (55.56), *maa://30050 (51.72)"
$ws.Range("AF29").Value = "*maa://24080 (69.05), ***maa://34960 (8.33), *maa://42865 (77.14)"
$ws.Range("AB30").Value = "maa://42979 (96.43), maa://45045 (100.0)"
$ws.Range("AF38").Value = "maa://36697 (86.26)"
$ws.Range("H62").Value = "maa://42981 (96.55), maa://43903 (100.0)"
